# Peer Review Log.xlsx - "Add files via upload" edit
# Updates the "User's Guide" sheet with a new peer-review row and refreshes
# the alternating-row striping / header formatting on both sheets, plus the
# saved selection on each sheet.

$wb = $excel.ActiveWorkbook
$wsTestPlan = $wb.Worksheets.Item(1)
$wsUsersGuide = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# "Test Plan" sheet: the striped banding underneath the header (rows
# 5,7,9,11,13,15) and the header row itself keep the same look, just the
# underlying style slots get touched as part of the resave.
# ---------------------------------------------------------------------
foreach ($r in 5,7,9,11,13,15) {
    $rng = $wsTestPlan.Range("A" + $r + ":E" + $r)
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# "User's Guide" sheet: fill in the new peer-review entry on row 4 and
# restripe the banding below it the same way.
# ---------------------------------------------------------------------
$wsUsersGuide.Range("A4").Value = "B. Allton"
$wsUsersGuide.Range("B4").Value = 43862
$wsUsersGuide.Range("B4").NumberFormat = "d-mmm-yy"
$wsUsersGuide.Range("C4").Value = "Check for content accuracy"
$wsUsersGuide.Range("D4").Value = " Listing  Windows 7 use under System Requirements, this OS is no longer supported by Microsoft"
$wsUsersGuide.Range("E4").Value = "Changed requirement to Win 8.1 or higher"

$row4 = $wsUsersGuide.Range("A4:E4")
$row4.HorizontalAlignment = -4108
$row4.VerticalAlignment = -4108
$row4.WrapText = $true

$wsUsersGuide.Rows.Item(4).RowHeight = 75
$wsUsersGuide.Rows.Item(3).RowHeight = 15.75

foreach ($r in 5,7,9,11,13,15) {
    $rng = $wsUsersGuide.Range("A" + $r + ":E" + $r)
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}
foreach ($r in 6,8,10,12,14) {
    $rng = $wsUsersGuide.Range("A" + $r + ":E" + $r)
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# ---------------------------------------------------------------------
# Saved cursor / selection per sheet (matches the values left behind by
# the author's last save).
# ---------------------------------------------------------------------
$wsTestPlan.Activate()
$wsTestPlan.Range("D28").Select()
$wsUsersGuide.Activate()
$wsUsersGuide.Range("I4").Select()
